# "Generate Report for Handback"
#
# The localization-status report records, per source file / per target
# language, the outcome of the latest handback attempt. This handback run
# failed for 4baf1cf7-c2aa-4f32-ad58-b4aa1a325048 (row 3) in both the
# zh-cn and de-de targets: the returned file name didn't match the
# handoff file name that was sent out. Update the Status columns on all
# three sheets and record the failure detail (with the Error Detail
# column widened so the message is readable).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: per-language status columns (E = zh-cn, F = de-de) for
# the 4baf1cf7... row.
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-language detail sheets: Status column (C) for the same row.
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Error Detail column (P) gets the handback mismatch explanation, and the
# column is widened so the long message is legible.
$wsZhCn.Range("P3").Value = "Handback file name: ct5dmczh.fea is different with handoff file name: 4baf1cf7-c2aa-4f32-ad58-b4aa1a325048.e43a2e1a65f2495179f65726bf57bf6ef5fd3692.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

$wsDeDe.Range("P3").Value = "Handback file name: ct5dmczh.fea is different with handoff file name: 4baf1cf7-c2aa-4f32-ad58-b4aa1a325048.e43a2e1a65f2495179f65726bf57bf6ef5fd3692.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
